# Adds three new slides to the deck:
#   - "Demo" (title-only content, blank body)
#   - "Publish your module to the web"
#   - "Materials" (with hyperlinks to external sources)
#
# PowerPoint assigns new slide IDs in creation order (261, 262, 263) but the
# final deck displays them as 262, 263, 261 - i.e. the "Demo" slide was
# created first but ends up last. We reproduce that by creating the slides
# in id order, then moving the first one ("Demo") to the end.

$p = $ppt.ActivePresentation

$layoutTitleAndContent = 2   # ppLayoutText - matches "Title and Content" used by the other content slides

$sDemo = $p.Slides.Add($p.Slides.Count + 1, $layoutTitleAndContent)
$sPublish = $p.Slides.Add($p.Slides.Count + 1, $layoutTitleAndContent)
$sMaterials = $p.Slides.Add($p.Slides.Count + 1, $layoutTitleAndContent)

# ---------------------------------------------------------------------
# Slide: Demo
# ---------------------------------------------------------------------
$titleRange = $sDemo.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Demo"
$titleRange.LanguageID = "nl-NL"

# Content placeholder is left empty (matches the target: an empty paragraph,
# no visible text).

# ---------------------------------------------------------------------
# Slide: Publish your module to the web
# ---------------------------------------------------------------------
$titleRange = $sPublish.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Publish your module to the web"
$titleRange.LanguageID = "nl-NL"

$bodyRange = $sPublish.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Easy but public`rPublish your module to the PowerShell Gallery `r`rMore effort but private`rHost your PowerShell module on a private hosting platform (artifactory, Azure DevOps Artifact, file share etc)`rRegister the hosting platform as a PSRepository`r`r"

$bodyRange.Paragraphs(1,1).ParagraphFormat.Bullet.Visible = $false
$bodyRange.Paragraphs(4,1).ParagraphFormat.Bullet.Visible = $false
$bodyRange.Paragraphs(7,1).ParagraphFormat.Bullet.Visible = $false

$bodyRange.Paragraphs(1,1).LanguageID = "nl-NL"

# ---------------------------------------------------------------------
# Slide: Materials
# ---------------------------------------------------------------------
$titleRange = $sMaterials.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Materials"
$titleRange.LanguageID = "nl-NL"

$bodyRange = $sMaterials.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Scripts and presentation`rGitHub`r`rSources`rPowerShell on docs.microsoft.com`rPowerShellExplained.com`rSimpleTalk.com"

$bodyRange.Paragraphs(1,1).ParagraphFormat.Bullet.Visible = $false
$bodyRange.Paragraphs(4,1).ParagraphFormat.Bullet.Visible = $false

$bodyRange.Paragraphs(2,1).ActionSettings.Item(1).Hyperlink.Address = "https://github.com/sanderrossel/creating-powershell-modules"
$bodyRange.Paragraphs(5,1).ActionSettings.Item(1).Hyperlink.Address = "https://docs.microsoft.com/en-us/powershell/"
$bodyRange.Paragraphs(6,1).ActionSettings.Item(1).Hyperlink.Address = "https://powershellexplained.com/"
$bodyRange.Paragraphs(7,1).ActionSettings.Item(1).Hyperlink.Address = "https://www.red-gate.com/simple-talk/sysadmin/powershell/"

$bodyRange.Paragraphs(1,1).LanguageID = "nl-NL"

# ---------------------------------------------------------------------
# Reorder: Demo slide (created first, id 261) moves to the very end so the
# final order is ..., Publish (262), Materials (263), Demo (261).
# ---------------------------------------------------------------------
$sDemo.MoveTo($p.Slides.Count)
